$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "Datos actualizados" timestamp string in cell A1
$ws.Range("A1").Value = "Datos actualizados a 20 de Abril de 2020 a las 18:22"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 770076
$ws.Range("C4").Value = 5440
$ws.Range("D4").Value = 71489
$ws.Range("E4").Value = 657271
$ws.Range("G4").Value = 741
$ws.Range("H4").Value = 41316

# Row 13 - Rusia
$ws.Range("F13").Value = 700

# Row 15 - Brasil
$ws.Range("B15").Value = 39384
$ws.Range("C15").Value = 730
$ws.Range("E15").Value = 14750
$ws.Range("G15").Value = 42
$ws.Range("H15").Value = 2504

# Row 16 - Canada
$ws.Range("B16").Value = 35708
$ws.Range("C16").Value = 652
$ws.Range("D16").Value = 12197
$ws.Range("E16").Value = 21893

# Row 18 - Suiza
$ws.Range("E18").Value = 8717
$ws.Range("G18").Value = 34
$ws.Range("H18").Value = 1427

# Row 52 - Luxemburgo
$ws.Range("B52").Value = 3558
$ws.Range("C52").Value = 8
$ws.Range("D52").Value = 637
$ws.Range("E52").Value = 2846
$ws.Range("F52").Value = 32
$ws.Range("G52").Value = 2
$ws.Range("H52").Value = 75

# Row 59 - Argelia
$ws.Range("B59").Value = 2718
$ws.Range("C59").Value = 89
$ws.Range("D59").Value = 1099
$ws.Range("E59").Value = 1235
$ws.Range("G59").Value = 9
$ws.Range("H59").Value = 384

# Row 67 - Islandia
$ws.Range("E67").Value = 401
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = 10

# Row 82 - Afganistan
$ws.Range("D82").Value = 135
$ws.Range("E82").Value = 855
$ws.Range("G82").Value = 3
$ws.Range("H82").Value = 36
